$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per the Jan 14 2023 GitHub Actions refresh.
# Each cell keeps its original plain-text representation (prices/percentages are
# stored as literal strings, not numbers), so we force Text format before writing
# and then clear the format again so no stray style is left behind on the cell.
$cellUpdates = @(
    @{ Cell = "D2"; Value = "304.35" }
    @{ Cell = "E2"; Value = "6.52%" }
    @{ Cell = "D3"; Value = "31.84" }
    @{ Cell = "E3"; Value = "8.47%" }
    @{ Cell = "D4"; Value = "5.262" }
    @{ Cell = "E4"; Value = "3.93%" }
    @{ Cell = "D5"; Value = "0.07535" }
    @{ Cell = "E5"; Value = "12.24%" }
    @{ Cell = "D6"; Value = "7.820" }
    @{ Cell = "E6"; Value = "6.89%" }
    @{ Cell = "D7"; Value = "3.753" }
    @{ Cell = "E7"; Value = "9.07%" }
    @{ Cell = "D8"; Value = "1.481" }
    @{ Cell = "E8"; Value = "6.82%" }
    @{ Cell = "D9"; Value = "0.9139" }
    @{ Cell = "E9"; Value = "1.57%" }
    @{ Cell = "D10"; Value = "0.01660" }
    @{ Cell = "E10"; Value = "2,469.12%" }
    @{ Cell = "D11"; Value = "0.1700" }
    @{ Cell = "E11"; Value = "6.38%" }
    @{ Cell = "D12"; Value = "0.07571" }
    @{ Cell = "E12"; Value = "6.71%" }
    @{ Cell = "D13"; Value = "0.08073" }
    @{ Cell = "E13"; Value = "5.64%" }
    @{ Cell = "D14"; Value = "0.02985" }
    @{ Cell = "E14"; Value = "2.06%" }
    @{ Cell = "D15"; Value = "0.09899" }
    @{ Cell = "D16"; Value = "0.001500" }
    @{ Cell = "E16"; Value = "-5.14%" }
    @{ Cell = "D17"; Value = "0.04553" }
    @{ Cell = "E17"; Value = "1.19%" }
    @{ Cell = "D18"; Value = "0.006113" }
    @{ Cell = "E18"; Value = "-1.28%" }
    @{ Cell = "D19"; Value = "3.491" }
    @{ Cell = "E19"; Value = "1.22%" }
    @{ Cell = "D20"; Value = "2.233" }
    @{ Cell = "E20"; Value = "0.11%" }
    @{ Cell = "E21"; Value = "2.42%" }
    @{ Cell = "E22"; Value = "1.77%" }
    @{ Cell = "D23"; Value = "4.485" }
    @{ Cell = "E23"; Value = "15.00%" }
    @{ Cell = "D24"; Value = "0.1628" }
    @{ Cell = "E24"; Value = "4.46%" }
    @{ Cell = "E25"; Value = "0.97%" }
    @{ Cell = "D26"; Value = "0.004430" }
    @{ Cell = "E26"; Value = "1.41%" }
    @{ Cell = "E27"; Value = "19.48%" }
    @{ Cell = "D28"; Value = "0.0001738" }
    @{ Cell = "E28"; Value = "7.45%" }
    @{ Cell = "D40"; Value = "0.04512" }
    @{ Cell = "E40"; Value = "6.26%" }
    @{ Cell = "D41"; Value = "0.007218" }
    @{ Cell = "E41"; Value = "5.84%" }
    @{ Cell = "D42"; Value = "0.1345" }
    @{ Cell = "E42"; Value = "8.61%" }
    @{ Cell = "D43"; Value = "0.002246" }
    @{ Cell = "E43"; Value = "0.75%" }
    @{ Cell = "E44"; Value = "1.76%" }
    @{ Cell = "E45"; Value = "7.96%" }
    @{ Cell = "D46"; Value = "0.7091" }
    @{ Cell = "E46"; Value = "-63.94%" }
    @{ Cell = "E47"; Value = "-13.54%" }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    $range.NumberFormat = "@"
    $range.Value = $update.Value
    $range.ClearFormats()
}
